$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 331, shifting existing rows 331-421 down to 332-422.
$ws.Rows.Item(331).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(331, 1).Value = 3
$ws.Cells.Item(331, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(331, 3).Value = "Coquimbo"
$ws.Cells.Item(331, 4).Value = 44841
$ws.Cells.Item(331, 5).Value = 5
$ws.Cells.Item(331, 6).Value = 100112009
$ws.Cells.Item(331, 7).Value = "Acelga"
$ws.Cells.Item(331, 8).Value = "Sin especificar"
$ws.Cells.Item(331, 9).Value = "Primera"
$ws.Cells.Item(331, 10).Value = 280
$ws.Cells.Item(331, 11).Value = 2300
$ws.Cells.Item(331, 12).Value = 2500
$ws.Cells.Item(331, 13).Value = 2386
$ws.Cells.Item(331, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(331, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(331, 16).Value = 398
$ws.Cells.Item(331, 17).Value = 6
$ws.Cells.Item(331, 18).Value = "Hortaliza"
